$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 24,14
$arr[0,0] = 19.49967933333333
$arr[0,1] = 58.499038
$arr[0,2] = 0.01453409262904611
$arr[0,3] = 0.01453409262904611
$arr[0,4] = 3.0
$arr[0,5] = 1.0
$arr[0,6] = 0.9442423333333334
$arr[0,7] = 2.832727
$arr[0,8] = 0.006848500623481535
$arr[0,9] = 0.006848500623481536
$arr[0,10] = 18.41242271295845
$arr[0,11] = 165.711804416626
$arr[0,12] = 0.00009953674243176066
$arr[0,13] = 0.00009953674243176067
$arr[1,0] = 19.49967933333333
$arr[1,1] = 58.499038
$arr[1,2] = 0.01453409262904611
$arr[1,3] = 0.01453409262904611
$arr[1,4] = 3.0
$arr[1,5] = 1.0
$arr[1,6] = 82.477727
$arr[1,7] = 247.433181
$arr[1,8] = 0.5982031781913751
$arr[1,9] = 0.5982031781913751
$arr[1,10] = 1608.289228642209
$arr[1,11] = 14474.60305777988
$arr[1,12] = 0.008694340402823221
$arr[1,13] = 0.008694340402823221
$arr[2,0] = 19.49967933333333
$arr[2,1] = 58.499038
$arr[2,2] = 0.01453409262904611
$arr[2,3] = 0.01453409262904611
$arr[2,4] = 3.0
$arr[2,5] = 1.0
$arr[2,6] = 9.766934000000001
$arr[2,7] = 29.300802
$arr[2,8] = 0.07083865150630789
$arr[2,9] = 0.07083865150630789
$arr[2,10] = 190.4520810698307
$arr[2,11] = 1714.068729628476
$arr[2,12] = 0.001029575522709396
$arr[2,13] = 0.001029575522709396
$arr[3,0] = 19.49967933333333
$arr[3,1] = 58.499038
$arr[3,2] = 0.01453409262904611
$arr[3,3] = 0.01453409262904611
$arr[3,4] = 3.0
$arr[3,5] = 1.0
$arr[3,6] = 44.68687199999999
$arr[3,7] = 134.060616
$arr[3,8] = 0.3241096696788354
$arr[3,9] = 0.3241096696788355
$arr[3,10] = 871.3796744097118
$arr[3,11] = 7842.417069687407
$arr[3,12] = 0.004710639961081731
$arr[3,13] = 0.004710639961081732
$arr[4,0] = 378.538676
$arr[4,1] = 1135.616028
$arr[4,2] = 0.2821439310161206
$arr[4,3] = 0.2821439310161206
$arr[4,4] = 3.0
$arr[4,5] = 1.0
$arr[4,6] = 0.9442423333333334
$arr[4,7] = 2.832727
$arr[4,8] = 0.006848500623481535
$arr[4,9] = 0.006848500623481536
$arr[4,10] = 357.4322426831507
$arr[4,11] = 3216.890184148356
$arr[4,12] = 0.001932262887475433
$arr[4,13] = 0.001932262887475434
$arr[5,0] = 378.538676
$arr[5,1] = 1135.616028
$arr[5,2] = 0.2821439310161206
$arr[5,3] = 0.2821439310161206
$arr[5,4] = 3.0
$arr[5,5] = 1.0
$arr[5,6] = 82.477727
$arr[5,7] = 247.433181
$arr[5,8] = 0.5982031781913751
$arr[5,9] = 0.5982031781913751
$arr[5,10] = 31221.00957806945
$arr[5,11] = 280989.086202625
$arr[5,12] = 0.1687793962412514
$arr[5,13] = 0.1687793962412515
$arr[6,0] = 378.538676
$arr[6,1] = 1135.616028
$arr[6,2] = 0.2821439310161206
$arr[6,3] = 0.2821439310161206
$arr[6,4] = 3.0
$arr[6,5] = 1.0
$arr[6,6] = 9.766934000000001
$arr[6,7] = 29.300802
$arr[6,8] = 0.07083865150630789
$arr[6,9] = 0.07083865150630789
$arr[6,10] = 3697.162264939384
$arr[6,11] = 33274.46038445445
$arr[6,12] = 0.01998669560387074
$arr[6,13] = 0.01998669560387074
$arr[7,0] = 378.538676
$arr[7,1] = 1135.616028
$arr[7,2] = 0.2821439310161206
$arr[7,3] = 0.2821439310161206
$arr[7,4] = 3.0
$arr[7,5] = 1.0
$arr[7,6] = 44.68687199999999
$arr[7,7] = 134.060616
$arr[7,8] = 0.3241096696788354
$arr[7,9] = 0.3241096696788355
$arr[7,10] = 16915.70936146147
$arr[7,11] = 152241.3842531532
$arr[7,12] = 0.09144557628352297
$arr[7,13] = 0.091445576283523
$arr[8,0] = 481.5587156666667
$arr[8,1] = 1444.676147
$arr[8,2] = 0.3589299526510408
$arr[8,3] = 0.3589299526510408
$arr[8,4] = 3.0
$arr[8,5] = 1.0
$arr[8,6] = 0.9442423333333334
$arr[8,7] = 2.832727
$arr[8,8] = 0.006848500623481535
$arr[8,9] = 0.006848500623481536
$arr[8,10] = 454.7081253180967
$arr[8,11] = 4092.37312786287
$arr[8,12] = 0.00245813200451685
$arr[8,13] = 0.002458132004516851
$arr[9,0] = 481.5587156666667
$arr[9,1] = 1444.676147
$arr[9,2] = 0.3589299526510408
$arr[9,3] = 0.3589299526510408
$arr[9,4] = 3.0
$arr[9,5] = 1.0
$arr[9,6] = 82.477727
$arr[9,7] = 247.433181
$arr[9,8] = 0.5982031781913751
$arr[9,9] = 0.5982031781913751
$arr[9,10] = 39717.86828522596
$arr[9,11] = 357460.8145670336
$arr[9,12] = 0.2147130384239324
$arr[9,13] = 0.2147130384239324
$arr[10,0] = 481.5587156666667
$arr[10,1] = 1444.676147
$arr[10,2] = 0.3589299526510408
$arr[10,3] = 0.3589299526510408
$arr[10,4] = 3.0
$arr[10,5] = 1.0
$arr[10,6] = 9.766934000000001
$arr[10,7] = 29.300802
$arr[10,8] = 0.07083865150630789
$arr[10,9] = 0.07083865150630789
$arr[10,10] = 4703.352193041101
$arr[10,11] = 42330.1697373699
$arr[10,12] = 0.02542611383102267
$arr[10,13] = 0.02542611383102267
$arr[11,0] = 481.5587156666667
$arr[11,1] = 1444.676147
$arr[11,2] = 0.3589299526510408
$arr[11,3] = 0.3589299526510408
$arr[11,4] = 3.0
$arr[11,5] = 1.0
$arr[11,6] = 44.68687199999999
$arr[11,7] = 134.060616
$arr[11,8] = 0.3241096696788354
$arr[11,9] = 0.3241096696788355
$arr[11,10] = 21519.35268748073
$arr[11,11] = 193674.1741873265
$arr[11,12] = 0.1163326683915689
$arr[11,13] = 0.1163326683915689
$arr[12,0] = 10.909999
$arr[12,1] = 32.729997
$arr[12,2] = 0.008131771468556478
$arr[12,3] = 0.008131771468556478
$arr[12,4] = 3.0
$arr[12,5] = 1.0
$arr[12,6] = 0.9442423333333334
$arr[12,7] = 2.832727
$arr[12,8] = 0.006848500623481535
$arr[12,9] = 0.006848500623481536
$arr[12,10] = 10.30168291242433
$arr[12,11] = 92.715146211819
$arr[12,12] = 0.0000556904419724184
$arr[12,13] = 0.0000556904419724184
$arr[13,0] = 10.909999
$arr[13,1] = 32.729997
$arr[13,2] = 0.008131771468556478
$arr[13,3] = 0.008131771468556478
$arr[13,4] = 3.0
$arr[13,5] = 1.0
$arr[13,6] = 82.477727
$arr[13,7] = 247.433181
$arr[13,8] = 0.5982031781913751
$arr[13,9] = 0.5982031781913751
$arr[13,10] = 899.831919092273
$arr[13,11] = 8098.487271830456
$arr[13,12] = 0.004864451536816431
$arr[13,13] = 0.004864451536816431
$arr[14,0] = 10.909999
$arr[14,1] = 32.729997
$arr[14,2] = 0.008131771468556478
$arr[14,3] = 0.008131771468556478
$arr[14,4] = 3.0
$arr[14,5] = 1.0
$arr[14,6] = 9.766934000000001
$arr[14,7] = 29.300802
$arr[14,8] = 0.07083865150630789
$arr[14,9] = 0.07083865150630789
$arr[14,10] = 106.557240173066
$arr[14,11] = 959.015161557594
$arr[14,12] = 0.0005760437251900099
$arr[14,13] = 0.0005760437251900099
$arr[15,0] = 10.909999
$arr[15,1] = 32.729997
$arr[15,2] = 0.008131771468556478
$arr[15,3] = 0.008131771468556478
$arr[15,4] = 3.0
$arr[15,5] = 1.0
$arr[15,6] = 44.68687199999999
$arr[15,7] = 134.060616
$arr[15,8] = 0.3241096696788354
$arr[15,9] = 0.3241096696788355
$arr[15,10] = 487.5337288331279
$arr[15,11] = 4387.803559498151
$arr[15,12] = 0.002635585764577618
$arr[15,13] = 0.002635585764577619
$arr[16,0] = 98.48487833333333
$arr[16,1] = 295.454635
$arr[16,2] = 0.0734057375912918
$arr[16,3] = 0.07340573759129182
$arr[16,4] = 3.0
$arr[16,5] = 1.0
$arr[16,6] = 0.9442423333333334
$arr[16,7] = 2.832727
$arr[16,8] = 0.006848500623481535
$arr[16,9] = 0.006848500623481536
$arr[16,10] = 92.99359131551611
$arr[16,11] = 836.942321839645
$arr[16,12] = 0.0005027192396610838
$arr[16,13] = 0.000502719239661084
$arr[17,0] = 98.48487833333333
$arr[17,1] = 295.454635
$arr[17,2] = 0.0734057375912918
$arr[17,3] = 0.07340573759129182
$arr[17,4] = 3.0
$arr[17,5] = 1.0
$arr[17,6] = 82.477727
$arr[17,7] = 247.433181
$arr[17,8] = 0.5982031781913751
$arr[17,9] = 0.5982031781913751
$arr[17,10] = 8122.808908804881
$arr[17,11] = 73105.28017924393
$arr[17,12] = 0.04391154552459286
$arr[17,13] = 0.04391154552459287
$arr[18,0] = 98.48487833333333
$arr[18,1] = 295.454635
$arr[18,2] = 0.0734057375912918
$arr[18,3] = 0.07340573759129182
$arr[18,4] = 3.0
$arr[18,5] = 1.0
$arr[18,6] = 9.766934000000001
$arr[18,7] = 29.300802
$arr[18,8] = 0.07083865150630789
$arr[18,9] = 0.07083865150630789
$arr[18,10] = 961.8953066796967
$arr[18,11] = 8657.05776011727
$arr[18,12] = 0.005199963463793005
$arr[18,13] = 0.005199963463793006
$arr[19,0] = 98.48487833333333
$arr[19,1] = 295.454635
$arr[19,2] = 0.0734057375912918
$arr[19,3] = 0.07340573759129182
$arr[19,4] = 3.0
$arr[19,5] = 1.0
$arr[19,6] = 44.68687199999999
$arr[19,7] = 134.060616
$arr[19,8] = 0.3241096696788354
$arr[19,9] = 0.3241096696788355
$arr[19,10] = 4400.981152017239
$arr[19,11] = 39608.83036815516
$arr[19,12] = 0.02379150936324486
$arr[19,13] = 0.02379150936324487
$arr[20,0] = 352.659012
$arr[20,1] = 1057.977036
$arr[20,2] = 0.2628545146439442
$arr[20,3] = 0.2628545146439442
$arr[20,4] = 3.0
$arr[20,5] = 1.0
$arr[20,6] = 0.9442423333333334
$arr[20,7] = 2.832727
$arr[20,8] = 0.006848500623481535
$arr[20,9] = 0.006848500623481536
$arr[20,10] = 332.9955683619081
$arr[20,11] = 2996.960115257172
$arr[20,12] = 0.001800159307423988
$arr[20,13] = 0.001800159307423989
$arr[21,0] = 352.659012
$arr[21,1] = 1057.977036
$arr[21,2] = 0.2628545146439442
$arr[21,3] = 0.2628545146439442
$arr[21,4] = 3.0
$arr[21,5] = 1.0
$arr[21,6] = 82.477727
$arr[21,7] = 247.433181
$arr[21,8] = 0.5982031781913751
$arr[21,9] = 0.5982031781913751
$arr[21,10] = 29086.51371582573
$arr[21,11] = 261778.6234424315
$arr[21,12] = 0.1572404060619588
$arr[21,13] = 0.1572404060619588
$arr[22,0] = 352.659012
$arr[22,1] = 1057.977036
$arr[22,2] = 0.2628545146439442
$arr[22,3] = 0.2628545146439442
$arr[22,4] = 3.0
$arr[22,5] = 1.0
$arr[22,6] = 9.766934000000001
$arr[22,7] = 29.300802
$arr[22,8] = 0.07083865150630789
$arr[22,9] = 0.07083865150630789
$arr[22,10] = 3444.397294709208
$arr[22,11] = 30999.57565238287
$arr[22,12] = 0.01862025935972207
$arr[22,13] = 0.01862025935972207
$arr[23,0] = 352.659012
$arr[23,1] = 1057.977036
$arr[23,2] = 0.2628545146439442
$arr[23,3] = 0.2628545146439442
$arr[23,4] = 3.0
$arr[23,5] = 1.0
$arr[23,6] = 44.68687199999999
$arr[23,7] = 134.060616
$arr[23,8] = 0.3241096696788354
$arr[23,9] = 0.3241096696788355
$arr[23,10] = 15759.22812889046
$arr[23,11] = 141833.0531600142
$arr[23,12] = 0.08519368991483937
$arr[23,13] = 0.08519368991483939
$ws.Range("G2:T25").Value = $arr
Write-Host "Updated G2:T25 with new TPM values"
